# Adds two new departure rows (26 & 27) to the "Main Data" sheet for a new
# date, "Saturday, Jan 14" - data pulled in via the new internet-download panel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: flight #25 (Birmingham, FR3696) ---------------------------------
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(26, 3).Value = "9:45 AM"
$ws.Cells.Item(26, 4).Value = "FR3696"
$ws.Cells.Item(26, 5).Value = "Birmingham"
$ws.Cells.Item(26, 6).Value = "(BHX)"
$ws.Cells.Item(26, 7).Value = "Ryanair "
$ws.Cells.Item(26, 8).Value = "B738"
$ws.Cells.Item(26, 9).Value = "(EI-EVH)"
$ws.Cells.Item(26, 10).Value = "9:41 AM"
$ws.Cells.Item(26, 12).Value = "0 hours, -4 minutes"

# --- Row 27: flight #26 (London, FR2474) -------------------------------------
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(27, 3).Value = "9:55 AM"
$ws.Cells.Item(27, 4).Value = "FR2474"
$ws.Cells.Item(27, 5).Value = "London"
$ws.Cells.Item(27, 6).Value = "(STN)"
$ws.Cells.Item(27, 7).Value = "Ryanair "
$ws.Cells.Item(27, 8).Value = "B38M"
$ws.Cells.Item(27, 9).Value = "(EI-HMS)"
$ws.Cells.Item(27, 10).Value = "10:06 AM"
$ws.Cells.Item(27, 12).Value = "0 hours, 11 minutes"

# Columns K (11) and M (13) stay blank in every data row - copy an existing
# blank, default-styled cell into the new rows so the blank cells persist
# with the same style as the rest of the table instead of being dropped.
$ws.Range("K2").Copy($ws.Range("K26"))
$ws.Range("M2").Copy($ws.Range("M26"))
$ws.Range("K2").Copy($ws.Range("K27"))
$ws.Range("M2").Copy($ws.Range("M27"))
